$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = 3700
$ws.Range("C3").Value = 1900
$ws.Range("D3").Value = 2500
$ws.Range("D2").Select()
